$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "#FieldMonster[{}]"
$ws.Range("B1").Value = "name"

$ws.Range("B1").Select()
